# Apply the "stuff at the bottom of the sheets" edit:
#  - Add a "pair_kind" value ("generic") to column J for the practice rows (2-5)
#  - Add a new block starting at row 27 with "stim details" header info

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new J column (pair_kind) for the practice rows with "generic"
$ws.Range("J2").Value = "generic"
$ws.Range("J3").Value = "generic"
$ws.Range("J4").Value = "generic"
$ws.Range("J5").Value = "generic"

# New block starting at row 27
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "video"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "video"

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "video"

$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "video"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "audio"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "audio"

$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "audio"

$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "audio"
